# "Flash Packets" sheet is tabSelected, i.e. the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new END_OF_PAGE packet-type row (row 8) -----------------------
# Referenced before "id" below so it lands first in the shared-string table
# (matches shared-string order END_OF_PAGE=52, id=53 in the target file).
$ws.Range("A8").Value = "END_OF_PAGE"
$ws.Range("B8").Value = 255
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = "ID"

# --- Insert a new "id" sub-field into the GYRO row (row 5) -----------------
# Shift the existing gxh/ghl/gyh/gyl/gzh/gzl values one column to the right
# (E5:J5 -> F5:K5), then put the new "id" label in the vacated E5 cell.
$ws.Range("K5").Value = $ws.Range("J5").Value2
$ws.Range("J5").Value = $ws.Range("I5").Value2
$ws.Range("I5").Value = $ws.Range("H5").Value2
$ws.Range("H5").Value = $ws.Range("G5").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("F5").Value = $ws.Range("E5").Value2
$ws.Range("E5").Value = "id"
